# Update the cryptos price list with refreshed figures (GitHub Actions run).
# For cells whose new text looks like a plain number (e.g. "102.50"), the
# cell is first coerced to Text format so Excel keeps it as a literal string
# (preserving trailing zeros / thousands-dot formatting) instead of silently
# converting it to a numeric value; the style is then restored to "Normal"
# immediately afterwards so no stray cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.512.54"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "3.018.82"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.543"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.62%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "3.494.73"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "3.022.50"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -15.03%  "
$ws.Range("D19").Value = "51.504.73"
$ws.Range("E19").Value = "  +1.58%  "
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("E22").Value = "  +1.84%  "
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("E26").Value = "  +4.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.172"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.38%  "
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "50.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0451"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.07%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.31%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.283"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.29%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.19%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.48%  "
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").Value = "2.028.17"
$ws.Range("E49").Value = "  +1.47%  "
$ws.Range("D50").Value = "3.316.27"
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("E51").Value = "  +2.15%  "
